$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Prepare the new last row (26) by copying the formatting (style) from row 24,
#    which carries the same "data row" style (s=2) as the row we are about to push down.
$ws.Range("A24:C24").Copy()
$ws.Range("A26:C26").PasteSpecial(-4122)

# 2) Update row 24: description cell (B) changes from "Testing Required" to "Done".
$ws.Cells.Item(24, 2).Value = "Done"

# 3) Insert the new row 25 content (TSID=ImportMojio, Description=Done, Runmode=N).
#    Row 25 already exists in the sheet, so writing to it preserves its existing style.
$ws.Cells.Item(25, 1).Value = "ImportMojio"
$ws.Cells.Item(25, 2).Value = "Done"
$ws.Cells.Item(25, 3).Value = "N"

# 4) Populate row 26 with the data that used to live in row 25
#    (TSID=ExportEvents, Description=<blank>, Runmode=Y).
$ws.Cells.Item(26, 1).Value = "ExportEvents"
$ws.Cells.Item(26, 2).Value = $null
$ws.Cells.Item(26, 3).Value = "Y"

# 5) Move the selection to match the saved view state.
$ws.Range("B19").Select()
